# Applies the update described in the commit:
# - Orders sheet: append 10 new flower rows (C52:C61)
# - Summary sheet: extend the zero-padded string in G2 by 10 more "0" characters

$wb = $excel.ActiveWorkbook

# --- Orders sheet: add new rows 52-61 in column C ---
$ordersSheet = $wb.Worksheets.Item("Orders")

$newFlowerRows = @(
    @(52, "517_鼠尾粉色_veronica pink_undefined_1bunch"),
    @(53, "396_米花 白_rice flower white_undefined_1bunch"),
    @(54, "490_米花 粉_rice flower pink_undefined_1bunch"),
    @(55, "484_天鹅绒_Star of Bethlehem_undefined_1bunch"),
    @(56, "594_绿毛球_undefined_undefined_1bunch"),
    @(57, "470_海芋白_Calla Lily_undefined_1bunch"),
    @(58, "745_海芋红_Calla Lily_undefined_1bunch"),
    @(59, "775_海芋黑_Calla Lily_undefined_1bunch"),
    @(60, "441_蓝星球_Echinops_undefined_1bunch"),
    @(61, "529_针垫_undefined_undefined_1bunch")
)

foreach ($item in $newFlowerRows) {
    $rowIndex = $item[0]
    $flowerName = $item[1]
    $ordersSheet.Cells.Item($rowIndex, 3).Value = $flowerName
}

# --- Summary sheet: extend G2 zero-padded code string ---
# The value is a long run of digits, so Excel's COM layer would normally
# coerce a plain .Value assignment into a (scientific-notation) Number.
# Force text storage via the "@" number format, then strip the temporary
# format back off so the cell keeps the workbook's default style, matching
# how the original text value was stored.
$summarySheet = $wb.Worksheets.Item("Summary")
$g2 = $summarySheet.Range("G2")
$g2.NumberFormat = "@"
$g2.Value = "05200000000000000000000000000000000000000000000000000000000000"
$g2.ClearFormats()

